$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Commercial Candida Rugosa lipase"
$ws.Range("A4").Value = "Biodegradation of Lipid-rich Waste Water by Combination of Microwave Irradiation and Lipase Immobilized on Chitosan"
$ws.Range("B4").Value = 2006
$ws.Range("D4").Value = "Free fatty acids (titrimetric)/Enzyme activity"

$ws.Range("A5").Value = "Effect of enzymatic pretreatment on the anaerobic digestion of milk fat for biogas production"
$ws.Range("B5").Value = 2015
$ws.Range("C5").Value = "Sigma Ladrich Candida Rugosa lipase"
$ws.Range("D5").Value = "Methane Production/COD/Free fatty acids"

$ws.Range("A4:D5").WrapText = $true
$ws.Range("A4:D5").HorizontalAlignment = -4131
$ws.Range("A4:D5").VerticalAlignment = -4160
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 45

$ws.Range("C11").Select()
